# Apply the "Compare" method change: new compare-method column (E) values
# plus a new threshold column (F) for a couple of rows, and refresh the
# window/selection view state to match the author's last interactive session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Window view state -------------------------------------------------
$excel.ActiveWindow.Left = 9720
$excel.ActiveWindow.Top = 220
$excel.ActiveWindow.Width = 25600
$excel.ActiveWindow.Height = 15520

# --- Data edits ----------------------------------------------------------
# Row 2: compare method -> VQSR, threshold 0.3
$ws.Range("E2").Value = "VQSR"
$ws.Range("F2").Value = 0.3

# Row 3: compare result -> PASS
$ws.Range("E3").Value = "PASS"

# Row 4: compare result -> FAIL
$ws.Range("E4").Value = "FAIL"

# Row 5: compare result -> PASS, threshold 0.004
$ws.Range("E5").Value = "PASS"
$ws.Range("F5").Value = 0.004

# Row 6: compare result -> PASS, threshold 0.0003
$ws.Range("E6").Value = "PASS"
$ws.Range("F6").Value = 0.0003

# --- Selection -------------------------------------------------------------
$ws.Range("F4").Select()
